$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: domain "https://www.nbg.gr/el/retail/housing-loans/Calculator/"
# now has a cookie present ("NBGPublicSite")
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "NBGPublicSite"

# Row 18: domain "https://microsites.nbg.gr/mobilebanking"
# no longer has the "NBGPUBLICConsent" cookie set
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = $null
